$d = $word.ActiveDocument

# Locate the final (last) paragraph in the document - an empty paragraph
# (w:pPr/w:jc=both) that sits right before the sectPr. The new TypeScript
# "Category" code block (plus a leading page break) is inserted immediately
# before it, and then this paragraph receives one new run of explanatory
# text.
$lastParagraph = $d.Paragraphs.Last
$insertionRange = $lastParagraph.Range
$insertionRange.Collapse(1)  # wdCollapseStart

$newParagraphsXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p>
      <w:r>
        <w:br w:type="page"/>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="1E1E1E"/>
        <w:spacing w:after="0" w:line="285" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="569CD6"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>class</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="4EC9B0"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>Category</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> {</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="1E1E1E"/>
        <w:spacing w:after="0" w:line="285" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="9CDCFE"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>id</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">?: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="4EC9B0"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>string</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>;</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="1E1E1E"/>
        <w:spacing w:after="0" w:line="285" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="9CDCFE"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>name</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="4EC9B0"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>string</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>;</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="1E1E1E"/>
        <w:spacing w:after="0" w:line="285" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="9CDCFE"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>description</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="4EC9B0"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>string</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>;</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="1E1E1E"/>
        <w:spacing w:after="0" w:line="285" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="9CDCFE"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>created_at</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="4EC9B0"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>Date</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>;</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="1E1E1E"/>
        <w:spacing w:after="0" w:line="285" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="1E1E1E"/>
        <w:spacing w:after="0" w:line="285" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="569CD6"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>constructor</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>() {</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="1E1E1E"/>
        <w:spacing w:after="0" w:line="285" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">        </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="C586C0"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>if</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>(!</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="569CD6"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>this</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="9CDCFE"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>id</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>) {</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="1E1E1E"/>
        <w:spacing w:after="0" w:line="285" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">            </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="569CD6"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>this</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="9CDCFE"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>id</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> = </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="DCDCAA"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>uuidv4</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>()</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="1E1E1E"/>
        <w:spacing w:after="0" w:line="285" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>        }</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="1E1E1E"/>
        <w:spacing w:after="0" w:line="285" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>    }</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="1E1E1E"/>
        <w:spacing w:after="0" w:line="285" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>}</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="1E1E1E"/>
        <w:spacing w:after="0" w:line="285" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="C586C0"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>export</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="C586C0"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>default</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
          <w:color w:val="4EC9B0"/>
          <w:kern w:val="0"/>
          <w:sz w:val="21"/>
          <w:szCs w:val="21"/>
          <w:lang w:eastAsia="pt-BR"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>Category</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertionRange.InsertXML($newParagraphsXml)

# Re-fetch the (now shifted) final paragraph and append the explanatory
# run describing the constructor method that was added.
$finalParagraph = $d.Paragraphs.Last
$finalParagraph.Range.InsertAfter("constructor(){} -> Adiciona m" + [char]0x00E9 + "todo " + [char]0x00E0 + " tipagem a ser exportada como model")

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
